# Update cryptos list (Coin / Link / Price / Volume(1h)) with the latest
# scraped figures. Price values in column D are stored as text (some
# contain multiple "." thousands separators), so a leading apostrophe is
# used wherever the new value would otherwise be auto-parsed as a number,
# to keep it as a text cell just like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "48.098.00"
$ws.Cells.Item(2, 5).Value = "  +0.68%  "
$ws.Cells.Item(3, 4).Value = "2.501.95"
$ws.Cells.Item(3, 5).Value = "  +0.19%  "
$ws.Cells.Item(4, 5).Value = "  -0.13%  "
$ws.Cells.Item(5, 4).Value = "'320.04"
$ws.Cells.Item(6, 4).Value = "'107.36"
$ws.Cells.Item(6, 5).Value = "  -1.37%  "
$ws.Cells.Item(7, 4).Value = "'0.526"
$ws.Cells.Item(7, 5).Value = "  +0.39%  "
$ws.Cells.Item(8, 5).Value = "  -0.09%  "
$ws.Cells.Item(9, 4).Value = "'0.541"
$ws.Cells.Item(9, 5).Value = "  -1.59%  "
$ws.Cells.Item(10, 5).Value = "  -0.92%  "
$ws.Cells.Item(11, 4).Value = "'20.12"
$ws.Cells.Item(11, 5).Value = "  +6.07%  "
$ws.Cells.Item(12, 5).Value = "  -0.30%  "
$ws.Cells.Item(13, 5).Value = "  +0.00%  "
$ws.Cells.Item(14, 4).Value = "'7.09"
$ws.Cells.Item(15, 4).Value = "2.893.20"
$ws.Cells.Item(15, 5).Value = "  +0.16%  "
$ws.Cells.Item(16, 4).Value = "2.503.28"
$ws.Cells.Item(16, 5).Value = "  +0.22%  "
$ws.Cells.Item(17, 4).Value = "'0.834"
$ws.Cells.Item(17, 5).Value = "  -1.82%  "
$ws.Cells.Item(18, 4).Value = "47.978.36"
$ws.Cells.Item(19, 4).Value = "'12.95"
$ws.Cells.Item(19, 5).Value = "  -1.59%  "
$ws.Cells.Item(20, 5).Value = "  +0.82%  "
$ws.Cells.Item(21, 5).Value = "  -0.17%  "
$ws.Cells.Item(22, 5).Value = "  -0.75%  "
$ws.Cells.Item(23, 4).Value = "'276.28"
$ws.Cells.Item(23, 5).Value = "  +11.56%  "
$ws.Cells.Item(24, 4).Value = "'71.54"
$ws.Cells.Item(24, 5).Value = "  +1.15%  "
$ws.Cells.Item(25, 5).Value = "  -1.27%  "
$ws.Cells.Item(26, 5).Value = "  -0.11%  "
$ws.Cells.Item(27, 4).Value = "'25.88"
$ws.Cells.Item(27, 5).Value = "  +0.04%  "
$ws.Cells.Item(28, 5).Value = "  +10.03%  "
$ws.Cells.Item(29, 2).Value = "Kaspa"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(29, 4).Value = "'0.141"
$ws.Cells.Item(29, 5).Value = "  +1.75%  "
$ws.Cells.Item(30, 2).Value = "Cosmos"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(30, 4).Value = "'9.72"
$ws.Cells.Item(30, 5).Value = "  -2.46%  "
$ws.Cells.Item(31, 4).Value = "'35.08"
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 4).Value = "'49.71"
$ws.Cells.Item(32, 5).Value = "  -0.39%  "
$ws.Cells.Item(33, 4).Value = "'19.48"
$ws.Cells.Item(33, 5).Value = "  -1.91%  "
$ws.Cells.Item(34, 5).Value = "  -0.18%  "
$ws.Cells.Item(35, 5).Value = "  -0.95%  "
$ws.Cells.Item(36, 4).Value = "'0.0782"
$ws.Cells.Item(36, 5).Value = "  -0.90%  "
$ws.Cells.Item(37, 5).Value = "  -0.90%  "
$ws.Cells.Item(38, 5).Value = "  -0.57%  "
$ws.Cells.Item(39, 4).Value = "'2.87"
$ws.Cells.Item(39, 5).Value = "  -2.77%  "
$ws.Cells.Item(40, 5).Value = "  -0.41%  "
$ws.Cells.Item(41, 4).Value = "'121.28"
$ws.Cells.Item(41, 5).Value = "  +1.53%  "
$ws.Cells.Item(42, 5).Value = "  -0.03%  "
$ws.Cells.Item(43, 4).Value = "'21.46"
$ws.Cells.Item(43, 5).Value = "  -4.44%  "
$ws.Cells.Item(44, 4).Value = "'0.0302"
$ws.Cells.Item(44, 5).Value = "  +1.66%  "
$ws.Cells.Item(45, 4).Value = "2.022.56"
$ws.Cells.Item(45, 5).Value = "  +1.01%  "
$ws.Cells.Item(46, 4).Value = "'3.13"
$ws.Cells.Item(46, 5).Value = "  +2.82%  "
$ws.Cells.Item(47, 5).Value = "  -1.62%  "
$ws.Cells.Item(48, 5).Value = "  +1.69%  "
$ws.Cells.Item(49, 4).Value = "'9.00"
$ws.Cells.Item(49, 5).Value = "  -0.08%  "
$ws.Cells.Item(50, 5).Value = "  +1.21%  "
$ws.Cells.Item(51, 4).Value = "'80.51"
$ws.Cells.Item(51, 5).Value = "  +3.60%  "
